$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price/Volume columns to Text format so numeric-looking
# values (e.g. "610.41") are stored as strings, matching the source data feed.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Update price and volume columns for rows 2-49 based on the latest crypto
# market snapshot
$ws.Range("D2").Value = '69.489.10'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '3.490.86'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '610.41'
$ws.Range("E5").Value = '  +4.73%  '
$ws.Range("D6").Value = '185.98'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.216'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '53.09'
$ws.Range("E11").Value = '  -2.52%  '
$ws.Range("D12").Value = '0.0000309'
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("D13").Value = '9.52'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '4.030.11'
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").Value = '604.51'
$ws.Range("E15").Value = '  +6.00%  '
$ws.Range("D16").Value = '69.452.21'
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("D17").Value = '12.66'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").Value = '18.87'
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").Value = '3.490.75'
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '0.988'
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("E22").Value = '  -2.58%  '
$ws.Range("D23").Value = '105.69'
$ws.Range("E23").Value = '  +11.38%  '
$ws.Range("D24").Value = '4.64'
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").Value = '5.05'
$ws.Range("E25").Value = '  +3.05%  '
$ws.Range("D26").Value = '3.03'
$ws.Range("E26").Value = '  +2.57%  '
$ws.Range("D27").Value = '10.95'
$ws.Range("E27").Value = '  -2.53%  '
$ws.Range("D28").Value = '9.77'
$ws.Range("E28").Value = '  +6.74%  '
$ws.Range("D29").Value = '33.68'
$ws.Range("E29").Value = '  +3.71%  '
$ws.Range("D30").Value = '6.99'
$ws.Range("E30").Value = '  -3.20%  '
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").Value = '4.00'
$ws.Range("E32").Value = '  +18.13%  '
$ws.Range("D33").Value = '0.116'
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").Value = '63.20'
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("E35").Value = '  -6.48%  '
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = '524.13'
$ws.Range("E37").Value = '  -3.97%  '
$ws.Range("D38").Value = '0.397'
$ws.Range("E38").Value = '  -3.83%  '
$ws.Range("D39").Value = '3.614.96'
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  +5.27%  '
$ws.Range("D41").Value = '36.75'
$ws.Range("E41").Value = '  -3.05%  '
$ws.Range("D42").Value = '0.0₃0777'
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").Value = '0.0459'
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("D45").Value = '2.95'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("E46").Value = '  +3.02%  '
$ws.Range("D47").Value = '3.33'
$ws.Range("E47").Value = '  -4.90%  '
$ws.Range("D48").Value = '8.82'
$ws.Range("E48").Value = '  -5.53%  '
$ws.Range("E49").Value = '  +0.37%  '

# Rows 50 and 51 swap coin identity (OceanProtocol and FLOKI swap rank order)
# and receive updated price / volume figures
$ws.Range("B50").Value = "OceanProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").Value = "  -9.18%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000242"
$ws.Range("E51").Value = "  -8.50%  "

# Restore the original (default) cell formatting so no visible style change
# is introduced by the temporary Text format above.
$dataRange.ClearFormats()
